# RaspberryPiList.xlsx update:
#  - AB-RPi02 status corrected from "Deployed, not configured" to "Running"
#  - Two new Service Block deployments appended to the RasPi inventory
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix AB-RPi02 status (row 9, column E)
$ws.Range("E9").Value = "Running"

# Insert two new rows for the new Service Block deployments after the
# existing last data row (row 15)
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(17).Insert()

$ws.Range("A16").Value = "FB-RPi01"
$ws.Range("B16").Value = "Service Block – Lower Ground Floor"
$ws.Range("C16").Value = "192.168.136.57"
$ws.Range("D16").Value = "b8:27:eb:73:d1:fd"
$ws.Range("E16").Value = "Running"

$ws.Range("A17").Value = "FB-RPi01"
$ws.Range("B17").Value = "Service Block – Upper Ground Floor"
$ws.Range("C17").Value = "192.168.136.59"
$ws.Range("D17").Value = "b8:27:eb:8d:2a:44"
$ws.Range("E17").Value = "Running"

# Keep the active selection on the last filled cell, matching the
# authored workbook's saved cursor position
$ws.Range("D17").Select()
